# Regenerate merged AHB files
# - Rename the diff-comparison header columns from generic "_old"/"_new"
#   suffixes to the concrete AHB format-version suffixes "_FV2304"/"_FV2310".
# - Turn the data range into a native Excel table ("Table1") with autofilter.
# - Freeze the header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (A1:U1) -----------------------------------------
$ws.Range("A1").Value = "Segmentname_FV2304"
$ws.Range("B1").Value = "Segmentgruppe_FV2304"
$ws.Range("C1").Value = "Segment_FV2304"
$ws.Range("D1").Value = "Datenelement_FV2304"
$ws.Range("E1").Value = "Segment ID_FV2304"
$ws.Range("F1").Value = "Code_FV2304"
$ws.Range("G1").Value = "Qualifier_FV2304"
$ws.Range("H1").Value = "Beschreibung_FV2304"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2304"
$ws.Range("J1").Value = "Bedingung_FV2304"
$ws.Range("K1").Value = "diff"
$ws.Range("L1").Value = "Segmentname_FV2310"
$ws.Range("M1").Value = "Segmentgruppe_FV2310"
$ws.Range("N1").Value = "Segment_FV2310"
$ws.Range("O1").Value = "Datenelement_FV2310"
$ws.Range("P1").Value = "Segment ID_FV2310"
$ws.Range("Q1").Value = "Code_FV2310"
$ws.Range("R1").Value = "Qualifier_FV2310"
$ws.Range("S1").Value = "Beschreibung_FV2310"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("U1").Value = "Bedingung_FV2310"

# --- 2. Turn A1:U92 into a table, "Table1", with the autofilter ------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U92"), 0, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (split below row 1) --------------------------
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
